$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "level" column (G) values for every question row (2..36 rows of data).
$levels = @(2,1,2,3,3,1,3,2,1,2,2,1,3,3,2,3,1,2,3,2,2,1,2,3,2,3,2,1,1,3,3,1,2,3,2,1)

for ($i = 0; $i -lt $levels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 7).Value = $levels[$i]
}

# Update the view: drop the old scroll position, zoom to 90%, and move the
# active selection from A19 to A18.
$excel.ActiveWindow.Zoom = 90
$ws.Range("A18").Select()
